# The deck ships two DrawingML themes:
#   ppt/theme/theme1.xml -> "Office Theme" (the stock Office palette)
#   ppt/theme/theme2.xml -> "Integral"     (the palette actually applied
#                                            to the one slide master, and
#                                            therefore to every slide)
# The authored edit swaps the two themes' content, so the deck's live
# design becomes the stock "Office Theme" palette. Re-create that by
# pushing the Office Theme's twelve scheme colours (dk1, lt1, dk2, lt2,
# accent1-6, hlink, folHlink - the standard PpThemeColorSchemeIndex
# order) onto the slide master's live ThemeColorScheme, which is the
# PowerPoint object model's handle onto the active theme part.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

$officeThemeColors = @{
    1  = 0          # dk1      000000
    2  = 16777215   # lt1      FFFFFF
    3  = 6968388     # dk2      44546A
    4  = 15132391    # lt2      E7E6E6
    5  = 13998939    # accent1  5B9BD5
    6  = 3243501      # accent2  ED7D31
    7  = 10855845     # accent3  A5A5A5
    8  = 49407        # accent4  FFC000
    9  = 12874308     # accent5  4472C4
    10 = 4697456      # accent6  70AD47
    11 = 12673797     # hlink    0563C1
    12 = 7491477      # folHlink 954F72
}

foreach ($index in 1..12) {
    $tcs.Item($index).RGB = $officeThemeColors[$index]
}
